$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to keep the literal text value instead of Excel
    # auto-converting numeric-looking strings (e.g. "291.09", "0.07644")
    # into floating point numbers, while leaving the cell's style
    # untouched (no persistent "Text" number format applied).
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "22.382.98"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.570.96"
$ws.Range("E3").Value = "  +0.25%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.00%  "

# Row 6 - BNB
Set-TextValue "D6" "291.09"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +2.06%  "

# Row 8 - OKB
Set-TextValue "D8" "49.94"
$ws.Range("E8").Value = "  +1.17%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3420"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.07644"
$ws.Range("E10").Value = "  +0.51%  "

# Row 11 - Polygon
Set-TextValue "D11" "1.149"
$ws.Range("E11").Value = "  -1.70%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  +0.01%  "

# Row 13 - Solana
Set-TextValue "D13" "21.16"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.009"
$ws.Range("E14").Value = "  -0.69%  "

# Row 15 - Chainlink
Set-TextValue "D15" "6.929"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.569.41"
$ws.Range("E16").Value = "  -0.53%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -0.57%  "

# Row 18 - Litecoin
$ws.Range("E18").Value = "  +0.87%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06737"
$ws.Range("E19").Value = "  -0.33%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.07%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  +1.52%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.203"
$ws.Range("E22").Value = "  -0.44%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -0.33%  "

# Row 24 - WrappedBTC
Set-TextValue "D24" "22.390.90"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +0.54%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.660"
$ws.Range("E26").Value = "  -10.86%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.19"
$ws.Range("E27").Value = "  +1.50%  "

# Row 28 - Monero
Set-TextValue "D28" "147.41"
$ws.Range("E28").Value = "  +1.24%  "

# Row 29 - HuobiToken
Set-TextValue "D29" "5.023"
$ws.Range("E29").Value = "  +1.47%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "126.15"
$ws.Range("E30").Value = "  +0.65%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-TextValue "D31" "1.744.21"
$ws.Range("E31").Value = "  -0.58%  "

# Row 32 - Filecoin
Set-TextValue "D32" "6.144"

# Row 33 - WEMIXTOKEN
Set-TextValue "D33" "2.007"
$ws.Range("E33").Value = "  +0.80%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "0.9810"
$ws.Range("E34").Value = "  -5.40%  "

# Row 35 - FraxShare
Set-TextValue "D35" "9.885"
$ws.Range("E35").Value = "  -3.97%  "

# Row 36 - Stellar
Set-TextValue "D36" "0.08479"
$ws.Range("E36").Value = "  +0.20%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.02544"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38 - TrustWalletToken
Set-TextValue "D38" "1.363"
$ws.Range("E38").Value = "  +9.56%  "

# Row 39 - Algorand
Set-TextValue "D39" "0.2317"
$ws.Range("E39").Value = "  -0.47%  "

# Row 40 - Hedera
Set-TextValue "D40" "0.06554"
$ws.Range("E40").Value = "  -0.10%  "

# Row 41 - InternetComputer(DFINITY)
Set-TextValue "D41" "5.414"
$ws.Range("E41").Value = "  -2.00%  "

# Row 42 - TheSandbox
Set-TextValue "D42" "0.6388"
$ws.Range("E42").Value = "  +0.24%  "

# Row 43 - Aptos
Set-TextValue "D43" "11.42"
$ws.Range("E43").Value = "  -3.52%  "

# Row 44 - was EnergySwap, now Frax (rows 44 and 45 swapped content)
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D44" "1.001"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - was Frax, now EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "14.03"
$ws.Range("E45").Value = "  -2.62%  "

# Row 46 - PancakeSwap
Set-TextValue "D46" "3.779"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47 - Decentraland
Set-TextValue "D47" "0.5971"
$ws.Range("E47").Value = "  -0.67%  "

# Row 48 - EOS
Set-TextValue "D48" "1.295"
$ws.Range("E48").Value = "  +1.88%  "

# Row 49 - NEARProtocol
Set-TextValue "D49" "2.087"
$ws.Range("E49").Value = "  -2.16%  "

# Row 50 - Quant
Set-TextValue "D50" "125.36"
$ws.Range("E50").Value = "  +1.36%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.07325"
$ws.Range("E51").Value = "  +0.51%  "
